# Edit script for Wellness workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename "Malik Boussaid" -> "Mehdi Boussaid" in the three rows that referenced it ---
$ws.Range("B765").Value = "Mehdi Boussaid"
$ws.Range("B778").Value = "Mehdi Boussaid"
$ws.Range("B792").Value = "Mehdi Boussaid"

# --- 2. Add the new rows 801-816 (training session on 2026-02-04, serial 46057) ---

# Copy formatting down first so new rows inherit the correct cell styles
$ws.Range("A800:I800").Copy()
$ws.Range("A801:I816").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Column G uses a different style depending on whether it holds text
# (s="1", like G798) or is left empty (s="2", like G800, already applied above)
$gTextRows = @(801, 802, 803, 804, 807, 808, 812, 814, 815, 816)
foreach ($gr in $gTextRows) {
    $ws.Range("G798").Copy()
    $ws.Range("G$gr").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

$nonBreak = [char]0x00A0

$rows = @(
    @{ R=801; B="Omar Benyounes";   D=4; E=3; F=1; G="Quadri$nonBreak";    H=6 },
    @{ R=802; B="Yoan Zouma";       D=5; E=4; F=2; G="Ischio$nonBreak";    H=5 },
    @{ R=803; B="Maé Clavel";       D=5; E=5; F=7; G="Genou/quadri";       H=2 },
    @{ R=804; B="Naim Ighbane";     D=3; E=6; F=3; G="Mollet$nonBreak";    H=2 },
    @{ R=805; B="Mehdi Boussaid";   D=5; E=0; F=0; G=$null;                H=7 },
    @{ R=806; B="Kamal Bafounta";   D=5; E=4; F=0; G=$null;                H=7 },
    @{ R=807; B="Romain Thunet";    D=7; E=6; F=5; G="Dos";                H=5 },
    @{ R=808; B="Nathanael Beta";   D=6; E=8; F=5; G="Dos";                H=7 },
    @{ R=809; B="Jeremie Laurent";  D=7; E=6; F=0; G=$null;                H=8 },
    @{ R=810; B="Ilan Ihaddadene";  D=5; E=5; F=0; G=$null;                H=3 },
    @{ R=811; B="Theo Owono";       D=6; E=2; F=0; G=$null;                H=5 },
    @{ R=812; B="Sofiane Belle";    D=5; E=4; F=4; G="Cheville$nonBreak";  H=7 },
    @{ R=813; B="Mattheo Haon";     D=5; E=3; F=0; G=$null;                H=4 },
    @{ R=814; B="Levy Ndoutoume";   D=5; E=5; F=1; G="Adducteur$nonBreak"; H=5 },
    @{ R=815; B="Hedi Nasri";       D=5; E=3; F=2; G="Ischio";             H=6 },
    @{ R=816; B="Naim Dhib";        D=6; E=4; F=3; G="Psoas";              H=4 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("A$r").Value = 46057
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = 70
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    if ($row.G) {
        $ws.Range("G$r").Value = $row.G
    } else {
        $ws.Range("G$r").ClearContents()
    }
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Formula = "=C$r*D$r"
}

# --- 3. Update the view: the window had scrolled further down and the
#        active selection moved from K795 to M805 ---
$ws.Range("M805").Select()
